$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" sheet: insert a new blank column before column N (14)
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N:N").Insert()

# Make this sheet the active one and set the new selection
$ws.Activate()
$ws.Range("R11").Select()
